# Auto-applied edits matching the target diff for cryptos.xlsx (Sheet1)
# Commit: "Updated cryptos list on Fri Sep 13 21:31:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cell updates (Coin name / Link / Volume% columns - never numeric-looking)
$textUpdates = @(
    @('E2', '  +3.89%  '),
    @('E3', '  +3.89%  '),
    @('E5', '  +3.31%  '),
    @('E6', '  +2.99%  '),
    @('E7', '  -0.09%  '),
    @('E8', '  +1.13%  '),
    @('E9', '  +5.49%  '),
    @('E10', '  +5.07%  '),
    @('E11', '  +2.53%  '),
    @('E12', '  -2.02%  '),
    @('E13', '  +5.24%  '),
    @('E14', '  +3.82%  '),
    @('E15', '  +3.69%  '),
    @('E16', '  +5.37%  '),
    @('E17', '  +2.76%  '),
    @('E18', '  +7.66%  '),
    @('E19', '  +4.33%  '),
    @('E20', '  +1.26%  '),
    @('E21', '  +2.40%  '),
    @('E22', '  +0.06%  '),
    @('E23', '  +3.36%  '),
    @('E24', '  +2.18%  '),
    @('E25', '  +1.12%  '),
    @('E26', '  +0.03%  '),
    @('E27', '  +0.87%  '),
    @('E28', '  +8.40%  '),
    @('E29', '  +4.02%  '),
    @('E30', '  -0.09%  '),
    @('E31', '  +3.36%  '),
    @('E32', '  +2.09%  '),
    @('E33', '  -0.84%  '),
    @('E35', '  +6.05%  '),
    @('E36', '  +1.60%  '),
    @('E37', '  +0.15%  '),
    @('E38', '  +0.54%  '),
    @('E39', '  +2.31%  '),
    @('E40', '  +10.58%  '),
    @('E41', '  +9.37%  '),
    @('E42', '  +2.82%  '),
    @('E43', '  -1.32%  '),
    @('E44', '  +1.94%  '),
    @('E45', '  +4.85%  '),
    @('E46', '  +1.37%  '),
    @('E47', '  +2.18%  '),
    @('B48', 'VeChain'),
    @('C48', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('E48', '  +3.54%  '),
    @('B49', 'Polygon'),
    @('C49', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @('E49', '  +5.07%  '),
    @('E50', '  -0.26%  '),
    @('E51', '  +5.51%  ')
)
foreach ($pair in $textUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Price column (D) updates: these are stored as TEXT in the workbook, but many
# of the new values look like plain numbers/decimals (e.g. "557.56", "0.400").
# Excel auto-converts a bare numeric-looking .Value assignment to a real number
# (losing formatting like trailing zeros, e.g. "0.400" -> 0.4). To preserve the
# original text semantics we temporarily force the cell to Text number format,
# assign the value, then restore General format + the default "Normal" style so
# no visible formatting/style change is left behind.
$priceUpdates = @(
    @('D2', '60.219.81'),
    @('D3', '2.435.42'),
    @('D5', '557.56'),
    @('D6', '139.37'),
    @('D9', '0.108'),
    @('D10', '5.83'),
    @('D11', '0.362'),
    @('D13', '24.99'),
    @('D14', '2.867.38'),
    @('D15', '60.076.34'),
    @('D17', '2.423.31'),
    @('D18', '11.49'),
    @('D19', '4.45'),
    @('D20', '335.54'),
    @('D23', '64.75'),
    @('D25', '8.57'),
    @('D27', '1.39'),
    @('D28', '0.0₃0797'),
    @('D29', '1.82'),
    @('D30', '171.28'),
    @('D32', '18.80'),
    @('D36', '4.28'),
    @('D39', '40.11'),
    @('D40', '0.417'),
    @('D41', '319.86'),
    @('D43', '143.14'),
    @('D44', '0.0965'),
    @('D45', '0.0527'),
    @('D47', '0.573'),
    @('D48', '0.0227'),
    @('D49', '0.400')
)
foreach ($pair in $priceUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
